$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 3684
$ws.Range("F5").Value = 3684
$ws.Range("F6").Value = 278
$ws.Range("F7").Value = 5210
$ws.Range("F8").Value = 559
$ws.Range("F9").Value = 390
$ws.Range("F10").Value = 215
$ws.Range("F11").Value = 716
$ws.Range("F13").Value = 116
$ws.Range("F14").Value = 40
$ws.Range("F15").Value = 718
$ws.Range("F16").Value = 331
$ws.Range("F17").Value = 41
$ws.Range("F18").Value = 94
$ws.Range("F19").Value = 162
$ws.Range("F22").Value = 5968
$ws.Range("F24").Value = 40
$ws.Range("F26").Value = 6200
$ws.Range("F27").Value = 21
$ws.Range("F29").Value = 3238
$ws.Range("F30").Value = 351
$ws.Range("F31").Value = 729
$ws.Range("F33").Value = 321
$ws.Range("F34").Value = 130
$ws.Range("F35").Value = 146
$ws.Range("F36").Value = 1087
$ws.Range("F37").Value = 92
$ws.Range("F38").Value = 26
$ws.Range("F40").Value = 899
$ws.Range("F41").Value = 1064
$ws.Range("F42").Value = 2041

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 1138

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 1138
$ws.Range("F7").Value = 3684
$ws.Range("F8").Value = 3684
$ws.Range("F9").Value = 278
$ws.Range("F10").Value = 5210
$ws.Range("F11").Value = 559
$ws.Range("F12").Value = 390
$ws.Range("F13").Value = 215
$ws.Range("F14").Value = 716
$ws.Range("F16").Value = 116
$ws.Range("F17").Value = 40
$ws.Range("F18").Value = 718
$ws.Range("F19").Value = 331
$ws.Range("F20").Value = 41
$ws.Range("F22").Value = 94
$ws.Range("F23").Value = 162
$ws.Range("F26").Value = 5968
$ws.Range("F28").Value = 40
$ws.Range("F30").Value = 6200
$ws.Range("F31").Value = 21
$ws.Range("F33").Value = 3238
$ws.Range("F34").Value = 351
$ws.Range("F35").Value = 729
$ws.Range("F37").Value = 321
$ws.Range("F39").Value = 130
$ws.Range("F40").Value = 146
$ws.Range("F41").Value = 1087
$ws.Range("F42").Value = 92
$ws.Range("F43").Value = 26
$ws.Range("F45").Value = 899
$ws.Range("F46").Value = 1064
$ws.Range("F48").Value = 2041
